# Refresh the cryptos price list (Coin / Link / Price / Volume(1h)) with the
# latest scraped figures. Most cells hold text (e.g. "52.130.20", "1.00")
# even though they look numeric, so a plain ".Value =" assignment would let
# Excel auto-convert them into real numbers. Set-TextValue works around
# that by briefly switching the cell to a text NumberFormat while the value
# is written, then restoring the "Normal" style so no stray style index is
# left behind in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: Bitcoin
$ws.Range('D2').Value = '52.130.20'
$ws.Range('E2').Value = '  +5.32%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.784.69'
$ws.Range('E3').Value = '  +5.69%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.03%  '

# Row 5: Solana
Set-TextValue $ws.Range('D5') '115.87'
$ws.Range('E5').Value = '  +3.21%  '

# Row 6: BNB
Set-TextValue $ws.Range('D6') '339.77'
$ws.Range('E6').Value = '  +4.21%  '

# Row 7: XRP
Set-TextValue $ws.Range('D7') '0.550'
$ws.Range('E7').Value = '  +5.15%  '

# Row 8: USDC
Set-TextValue $ws.Range('D8') '0.999'
$ws.Range('E8').Value = '  +0.00%  '

# Row 9: Cardano
Set-TextValue $ws.Range('D9') '0.578'
$ws.Range('E9').Value = '  +5.08%  '

# Row 10: Avalanche
Set-TextValue $ws.Range('D10') '41.97'
$ws.Range('E10').Value = '  +5.99%  '

# Row 11: Dogecoin
$ws.Range('E11').Value = '  +6.04%  '

# Row 12: Chainlink
Set-TextValue $ws.Range('D12') '20.09'
$ws.Range('E12').Value = '  +0.83%  '

# Row 13: TRON
$ws.Range('E13').Value = '  +2.25%  '

# Row 14: Polkadot
Set-TextValue $ws.Range('D14') '7.61'
$ws.Range('E14').Value = '  -0.05%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '3.224.76'
$ws.Range('E15').Value = '  +5.86%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '2.772.87'
$ws.Range('E16').Value = '  +5.68%  '

# Row 17: Polygon
Set-TextValue $ws.Range('D17') '0.883'
$ws.Range('E17').Value = '  +3.25%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '51.960.66'
$ws.Range('E18').Value = '  +5.09%  '

# Row 19: ImmutableX
Set-TextValue $ws.Range('D19') '3.21'
$ws.Range('E19').Value = '  +10.19%  '

# Row 20: InternetComputer(DFINITY)
$ws.Range('E20').Value = '  -0.83%  '

# Row 21: Uniswap
Set-TextValue $ws.Range('D21') '6.96'
$ws.Range('E21').Value = '  +4.54%  '

# Row 22: ShibaInu
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('E22').Value = '  +3.29%  '

# Row 23: BitcoinCash
Set-TextValue $ws.Range('D23') '278.16'
$ws.Range('E23').Value = '  +3.53%  '

# Row 24: Litecoin
Set-TextValue $ws.Range('D24') '70.12'
$ws.Range('E24').Value = '  +1.43%  '

# Row 25: PancakeSwap
Set-TextValue $ws.Range('D25') '2.74'
$ws.Range('E25').Value = '  +7.04%  '

# Row 26: EthereumClassic
$ws.Range('E26').Value = '  +2.71%  '

# Row 27: Dai
$ws.Range('E27').Value = '  -0.04%  '

# Row 28: Cosmos
Set-TextValue $ws.Range('D28') '10.22'
$ws.Range('E28').Value = '  +0.74%  '

# Row 29: Toncoin
$ws.Range('E29').Value = '  +1.07%  '

# Row 30: Kaspa
$ws.Range('E30').Value = '  +3.29%  '

# Row 31: InjectiveProtocol
$ws.Range('E31').Value = '  +0.94%  '

# Row 32: OKB
Set-TextValue $ws.Range('D32') '50.21'
$ws.Range('E32').Value = '  +1.17%  '

# Row 33: Filecoin
Set-TextValue $ws.Range('D33') '5.72'
$ws.Range('E33').Value = '  +4.56%  '

# Row 34: Hedera
Set-TextValue $ws.Range('D34') '0.0825'
$ws.Range('E34').Value = '  +1.37%  '

# Row 35: ARBITRUM -> FirstDigitalUSD
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D35') '1.00'
$ws.Range('E35').Value = '  -0.06%  '

# Row 36: FirstDigitalUSD -> ARBITRUM
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D36') '2.11'
$ws.Range('E36').Value = '  +4.32%  '

# Row 37: Celestia
Set-TextValue $ws.Range('D37') '18.92'
$ws.Range('E37').Value = '  -0.81%  '

# Row 38: RenderToken
Set-TextValue $ws.Range('D38') '4.97'
$ws.Range('E38').Value = '  +0.30%  '

# Row 39: LidoDAOToken
$ws.Range('E39').Value = '  +4.58%  '

# Row 40: VeChain
Set-TextValue $ws.Range('D40') '0.0377'
$ws.Range('E40').Value = '  +11.98%  '

# Row 41: Stacks
$ws.Range('E41').Value = '  +28.82%  '

# Row 42: WEMIXToken
$ws.Range('E42').Value = '  -1.13%  '

# Row 43: Stellar
Set-TextValue $ws.Range('D43') '0.115'
$ws.Range('E43').Value = '  +3.79%  '

# Row 44: EnergySwap
Set-TextValue $ws.Range('D44') '23.14'
$ws.Range('E44').Value = '  -0.49%  '

# Row 45: Monero
Set-TextValue $ws.Range('D45') '124.97'
$ws.Range('E45').Value = '  -3.56%  '

# Row 46: Maker
$ws.Range('D46').Value = '2.085.41'
$ws.Range('E46').Value = '  +1.09%  '

# Row 47: NEARProtocol
Set-TextValue $ws.Range('D47') '3.31'
$ws.Range('E47').Value = '  +0.21%  '

# Row 48: ApeXProtocol
$ws.Range('E48').Value = '  +3.46%  '

# Row 49: THORChain
Set-TextValue $ws.Range('D49') '5.56'
$ws.Range('E49').Value = '  +6.54%  '

# Row 50: FraxShare
Set-TextValue $ws.Range('D50') '8.96'
$ws.Range('E50').Value = '  +1.00%  '

# Row 51: SEI
Set-TextValue $ws.Range('D51') '0.890'
$ws.Range('E51').Value = '  +19.73%  '

Write-Host "Applied cryptos update."